$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New machine rows (A18, A19, A20 entered first, then A17 - matches the
# shared-string allocation order seen in the target workbook).
$ws.Range("A18").Value = "Butoniera"
$ws.Range("A19").Value = "Lucru manual"
$ws.Range("A20").Value = "Calcat."
$ws.Range("A17").Value = "Capse"

# Formatting for the new B column cells (B17:B20): left/center aligned,
# wrapped text, indented one level, explicit Calibri 11 font.
$rng = $ws.Range("B17:B20")
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4108
$rng.WrapText = $true
$rng.IndentLevel = 1
$rng.Font.Name = "Calibri"
$rng.Font.Size = 11

$ws.Range("D21").Select() | Out-Null
